# Auto-generated script to update cryptos price/volume columns
# to match the refreshed data snapshot described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each touched cell to plain Text before writing so Excel
# does not reinterpret values such as "1.040" or "6.820" as
# numbers (which would silently drop the meaningful trailing
# zero) - the source data are plain inline strings, not numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.058.87"
Set-TextValue $ws.Range("E2") "  -3.06%  "
Set-TextValue $ws.Range("D3") "1.715.61"
Set-TextValue $ws.Range("E3") "  -2.97%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "310.33"
Set-TextValue $ws.Range("E5") "  -5.64%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.07%  "
Set-TextValue $ws.Range("D7") "0.4798"
Set-TextValue $ws.Range("E7") "  +5.70%  "
Set-TextValue $ws.Range("D8") "0.3453"
Set-TextValue $ws.Range("E8") "  -2.18%  "
Set-TextValue $ws.Range("D9") "42.19"
Set-TextValue $ws.Range("E9") "  +0.68%  "
Set-TextValue $ws.Range("D10") "0.07248"
Set-TextValue $ws.Range("D11") "1.040"
Set-TextValue $ws.Range("E11") "  -4.97%  "
Set-TextValue $ws.Range("E12") "  -0.03%  "
Set-TextValue $ws.Range("D13") "19.74"
Set-TextValue $ws.Range("E13") "  -4.74%  "
Set-TextValue $ws.Range("D14") "5.839"
Set-TextValue $ws.Range("E14") "  -2.88%  "
Set-TextValue $ws.Range("D15") "1.717.08"
Set-TextValue $ws.Range("E15") "  -3.08%  "
Set-TextValue $ws.Range("D16") "6.820"
Set-TextValue $ws.Range("D17") "87.27"
Set-TextValue $ws.Range("E17") "  -5.80%  "
Set-TextValue $ws.Range("D18") "0.00001032"
Set-TextValue $ws.Range("E18") "  -2.45%  "
Set-TextValue $ws.Range("D19") "0.06378"
Set-TextValue $ws.Range("E19") "  -0.73%  "
Set-TextValue $ws.Range("E20") "  -0.04%  "
Set-TextValue $ws.Range("E21") "  -3.11%  "
Set-TextValue $ws.Range("D22") "5.626"
Set-TextValue $ws.Range("E22") "  -2.47%  "
Set-TextValue $ws.Range("D23") "27.109.46"
Set-TextValue $ws.Range("E23") "  -2.98%  "
Set-TextValue $ws.Range("E24") "  -4.16%  "
Set-TextValue $ws.Range("D25") "2.093"
Set-TextValue $ws.Range("E25") "  -0.35%  "
Set-TextValue $ws.Range("E26") "  -1.19%  "
Set-TextValue $ws.Range("D27") "150.71"
Set-TextValue $ws.Range("E27") "  -5.60%  "
Set-TextValue $ws.Range("D28") "1.910.97"
Set-TextValue $ws.Range("E28") "  -3.25%  "
Set-TextValue $ws.Range("D29") "2.062"
Set-TextValue $ws.Range("E29") "  -3.12%  "
Set-TextValue $ws.Range("D30") "120.70"
Set-TextValue $ws.Range("E30") "  -2.87%  "
Set-TextValue $ws.Range("D31") "1.038"
Set-TextValue $ws.Range("E31") "  -4.03%  "
Set-TextValue $ws.Range("D32") "0.09235"
Set-TextValue $ws.Range("E32") "  +0.51%  "
Set-TextValue $ws.Range("D33") "3.599"
Set-TextValue $ws.Range("E33") "  -1.83%  "
Set-TextValue $ws.Range("D34") "5.310"
Set-TextValue $ws.Range("E34") "  -5.28%  "
Set-TextValue $ws.Range("D35") "1.475"
Set-TextValue $ws.Range("E35") "  +6.95%  "
Set-TextValue $ws.Range("D36") "0.02180"
Set-TextValue $ws.Range("E36") "  -4.21%  "
Set-TextValue $ws.Range("D37") "0.05841"
Set-TextValue $ws.Range("E37") "  -4.13%  "
Set-TextValue $ws.Range("D38") "10.93"
Set-TextValue $ws.Range("E38") "  -7.50%  "
Set-TextValue $ws.Range("D39") "0.1983"
Set-TextValue $ws.Range("E39") "  -4.94%  "
Set-TextValue $ws.Range("D40") "1.000"
Set-TextValue $ws.Range("E40") "  -0.05%  "
Set-TextValue $ws.Range("D42") "0.5934"
Set-TextValue $ws.Range("E42") "  -5.00%  "
Set-TextValue $ws.Range("D43") "1.082"
Set-TextValue $ws.Range("E43") "  -8.12%  "
Set-TextValue $ws.Range("D44") "7.494"
Set-TextValue $ws.Range("E44") "  -3.95%  "
Set-TextValue $ws.Range("D45") "12.78"
Set-TextValue $ws.Range("E45") "  -3.56%  "
Set-TextValue $ws.Range("D46") "3.585"
Set-TextValue $ws.Range("E46") "  -4.06%  "
Set-TextValue $ws.Range("D47") "0.5562"
Set-TextValue $ws.Range("E47") "  -4.67%  "
Set-TextValue $ws.Range("D48") "118.60"
Set-TextValue $ws.Range("E48") "  -3.19%  "
Set-TextValue $ws.Range("E49") "  -5.63%  "
Set-TextValue $ws.Range("D50") "0.06636"
Set-TextValue $ws.Range("E50") "  -2.96%  "
Set-TextValue $ws.Range("D51") "1.087"
Set-TextValue $ws.Range("E51") "  -4.32%  "
